# Updated cryptos list with refreshed Price/Volume(1h) figures from the
# coinranking data pull, including the Maker/Bittensor rank swap (rows 38-39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.553.50"

$ws.Range("D3").Value = "3.495.38"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.49"
$ws.Range("E5").Value = "  +5.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.38"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -3.54%  "

$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.97"
$ws.Range("E11").Value = "  -3.19%  "

$ws.Range("E12").Value = "  -3.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.50"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").Value = "4.056.33"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "600.86"
$ws.Range("E15").Value = "  +5.03%  "

$ws.Range("D16").Value = "69.600.89"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.98"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.59"
$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").Value = "3.497.24"
$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.986"
$ws.Range("E21").Value = "  -1.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.25"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.23"
$ws.Range("E23").Value = "  +12.31%  "

$ws.Range("E24").Value = "  +4.23%  "

$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.05"
$ws.Range("E26").Value = "  +2.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -2.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.73"
$ws.Range("E28").Value = "  +5.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.35"
$ws.Range("E29").Value = "  +2.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.94"
$ws.Range("E30").Value = "  -3.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.14"
$ws.Range("E31").Value = "  +12.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.50"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.46"
$ws.Range("E34").Value = "  +0.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.17"
$ws.Range("E35").Value = "  -6.05%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.69"
$ws.Range("E37").Value = "  +7.19%  "

# Rows 38-39 swapped rank order: Maker now above Bittensor.
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.622.89"
$ws.Range("E38").Value = "  +1.34%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "512.54"
$ws.Range("E39").Value = "  -5.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.394"
$ws.Range("E40").Value = "  -4.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.80"
$ws.Range("E41").Value = "  -3.71%  "

$ws.Range("D42").Value = "0.0₃0775"
$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.92"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.141"
$ws.Range("E46").Value = "  +2.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  -4.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.76"
$ws.Range("E48").Value = "  -5.97%  "

$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.22"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("E51").Value = "  -8.77%  "
